$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): row 11
$ws.Range("B11").Value = "'48.28"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'13.63"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'61.91"
$ws.Range("D11").Style = "Normal"

# Employment (% of total): row 12
$ws.Range("B12").Value = "'14.66"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'33.13"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'47.79"
$ws.Range("D12").Style = "Normal"

# Enterprises (% of total): row 14
$ws.Range("B14").Value = "'75.76"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'21.38"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'97.14"
$ws.Range("D14").Style = "Normal"
